$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.517.89"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "3.847.80"
$ws.Range("E3").Value = "  +10.08%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'426.11"
$ws.Range("E5").Value = "  +8.96%  "
$ws.Range("D6").Value = "'130.89"
$ws.Range("E6").Value = "  +6.73%  "
$ws.Range("D7").Value = "3.842.20"
$ws.Range("E7").Value = "  +10.11%  "
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.729"
$ws.Range("E10").Value = "  +8.89%  "
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  +6.91%  "
$ws.Range("D12").Value = "'0.0000339"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "'41.41"
$ws.Range("E13").Value = "  +7.00%  "
$ws.Range("D14").Value = "'10.43"
$ws.Range("E14").Value = "  +13.52%  "
$ws.Range("D15").Value = "4.459.71"
$ws.Range("E15").Value = "  +11.65%  "
$ws.Range("D16").Value = "'15.89"
$ws.Range("E16").Value = "  +27.47%  "
$ws.Range("D17").Value = "3.874.98"
$ws.Range("E17").Value = "  +10.72%  "
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'20.00"
$ws.Range("E19").Value = "  +7.60%  "
$ws.Range("E20").Value = "  +7.72%  "
$ws.Range("D21").Value = "66.845.32"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").Value = "'413.93"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("D23").Value = "'15.04"
$ws.Range("E23").Value = "  +8.99%  "
$ws.Range("D24").Value = "'84.60"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("E25").Value = "  +7.96%  "
$ws.Range("D26").Value = "'37.67"
$ws.Range("E26").Value = "  +13.90%  "
$ws.Range("E27").Value = "  +13.74%  "
$ws.Range("E28").Value = "  +9.86%  "
$ws.Range("D29").Value = "'5.33"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Value = "'9.05"
$ws.Range("E30").Value = "  +34.19%  "
$ws.Range("D31").Value = "'718.26"
$ws.Range("E31").Value = "  +9.83%  "
$ws.Range("D32").Value = "'13.70"
$ws.Range("E32").Value = "  +15.80%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("E33").Value = "  +13.41%  "
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'39.13"
$ws.Range("E36").Value = "  +6.74%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.65"
$ws.Range("E38").Value = "  +40.36%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'55.70"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("E40").Value = "  +17.63%  "
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  +8.44%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +9.12%  "
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("E46").Value = "  +10.52%  "
$ws.Range("D47").Value = "'0.316"
$ws.Range("E47").Value = "  +15.17%  "
$ws.Range("D48").Value = "'141.98"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("E51").Value = "  +4.51%  "
